$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5
$ws.Range("C5").Value = "xcgbjaeitgbvzmt@gmail.com"
$ws.Range("D5").Value = "sjlrlLKJSB5"

# Update row 6
$ws.Range("C6").Value = "xiytnjkplizixwd@gmail.com"
$ws.Range("D6").Value = "yhoqqXQNPK5"

# Update row 7
$ws.Range("C7").Value = "wvadvqtseuzaftr@gmail.com"
$ws.Range("D7").Value = "glcbqYPPYZ5"

# Update selection to a single cell G20
$ws.Range("G20").Select()
